# Refresh market-price-derived columns (currentAveragePrice*, LevePrice*, LeveProfit*)
# on specific Leve rows across all 8 crafting-job sheets, per the scheduled market-data pull.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 33: Glazed and Confused
$ws.Range("H33").Value = 34483040
$ws.Range("I33").Value = 157.15
$ws.Range("K33").Value = 157.15
$ws.Range("M33").Value = 71.84999999999999
# Row 74: Adhesive of Antipathy
$ws.Range("H74").Value = 3100.4443
$ws.Range("I74").Value = 2500
$ws.Range("J74").Value = 3175.5
$ws.Range("K74").Value = 2500
$ws.Range("L74").Value = 3175.5
$ws.Range("M74").Value = -1564
$ws.Range("N74").Value = -5047.5
# Row 76: Warding Off Temptation
$ws.Range("H76").Value = 2986.6667
$ws.Range("I76").Value = 2985.7144
$ws.Range("J76").Value = 3000
$ws.Range("K76").Value = 2985.7144
$ws.Range("L76").Value = 3000
$ws.Range("M76").Value = -2670.7144
$ws.Range("N76").Value = -3630
# Row 77: It's Gonna Grow Back (L)
$ws.Range("H77").Value = 3100.4443
$ws.Range("I77").Value = 2500
$ws.Range("J77").Value = 3175.5
$ws.Range("K77").Value = 12500
$ws.Range("L77").Value = 15877.5
$ws.Range("M77").Value = -7820
$ws.Range("N77").Value = -25237.5
# Row 79: The Garden of Arcane Delights (L)
$ws.Range("H79").Value = 2986.6667
$ws.Range("I79").Value = 2985.7144
$ws.Range("J79").Value = 3000
$ws.Range("K79").Value = 2985.7144
$ws.Range("L79").Value = 3000
$ws.Range("M79").Value = -1893.7144
$ws.Range("N79").Value = -5184
# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 4443.25
$ws.Range("I132").Value = 1962.4706
$ws.Range("K132").Value = 5887.4118
$ws.Range("M132").Value = -3357.4118
# Row 135: For Tired Minds
$ws.Range("H135").Value = 33334002
$ws.Range("I135").Value = 454
$ws.Range("J135").Value = 100001096
$ws.Range("K135").Value = 4086
$ws.Range("L135").Value = 900009864
$ws.Range("M135").Value = -1551
$ws.Range("N135").Value = -900014934
# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 2138192.5
$ws.Range("I137").Value = 1457.5294
$ws.Range("J137").Value = 6174247.5
$ws.Range("K137").Value = 4372.5882
$ws.Range("L137").Value = 18522742.5
$ws.Range("M137").Value = -1822.5882
$ws.Range("N137").Value = -18527842.5
# Row 138: All-night Crafting
$ws.Range("H138").Value = 2600128.5
$ws.Range("I138").Value = 1447.7142
$ws.Range("J138").Value = 3574633.8
$ws.Range("K138").Value = 4343.142599999999
$ws.Range("L138").Value = 10723901.4
$ws.Range("M138").Value = 796.8574000000008
$ws.Range("N138").Value = -10734181.4

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 4607.8257
$ws.Range("I32").Value = 3038.3472
$ws.Range("K32").Value = 3038.3472
$ws.Range("M32").Value = -2751.3472
# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 1529.3778
$ws.Range("I61").Value = 1371.2354
$ws.Range("J61").Value = 2018.1818
$ws.Range("K61").Value = 1371.2354
$ws.Range("L61").Value = 2018.1818
$ws.Range("M61").Value = -1159.2354
$ws.Range("N61").Value = -2442.1818
# Row 63: Rivets Run through It
$ws.Range("H63").Value = 3287.5
$ws.Range("I63").Value = 2975
$ws.Range("J63").Value = 3350
$ws.Range("K63").Value = 2975
$ws.Range("L63").Value = 3350
$ws.Range("M63").Value = -2289
$ws.Range("N63").Value = -4722
# Row 66: A Riveting Revival (L)
$ws.Range("H66").Value = 3287.5
$ws.Range("I66").Value = 2975
$ws.Range("J66").Value = 3350
$ws.Range("K66").Value = 14875
$ws.Range("L66").Value = 16750
$ws.Range("M66").Value = -11443
$ws.Range("N66").Value = -23614
# Row 104: See Shields by the Sea Shore
$ws.Range("H104").Value = 28330
$ws.Range("J104").Value = 28330
$ws.Range("L104").Value = 28330
$ws.Range("N104").Value = -35318
# Row 133: Shielding My Students
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = $null
# Row 136: Metal with Mettle
$ws.Range("H136").Value = 1529.3778
$ws.Range("I136").Value = 1371.2354
$ws.Range("J136").Value = 2018.1818
$ws.Range("K136").Value = 4113.706200000001
$ws.Range("L136").Value = 6054.5454
$ws.Range("M136").Value = -1563.706200000001
$ws.Range("N136").Value = -11154.5454

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# Row 48: Colder than Steel
$ws.Range("H48").Value = 100000
$ws.Range("J48").Value = 100000
$ws.Range("L48").Value = 100000
$ws.Range("N48").Value = -100830
# Row 51: A Mixed Message
$ws.Range("H51").Value = 29649.5
$ws.Range("J51").Value = 29649.5
$ws.Range("L51").Value = 29649.5
$ws.Range("N51").Value = -30631.5
# Row 107: The Gold Experience
$ws.Range("H107").Value = 1058.9354
$ws.Range("I107").Value = 1074.2222
$ws.Range("J107").Value = 955.75
$ws.Range("K107").Value = 1074.2222
$ws.Range("L107").Value = 955.75
$ws.Range("M107").Value = 845.7778000000001
$ws.Range("N107").Value = -4795.75
# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 609516.4
$ws.Range("I134").Value = 1028556.2
$ws.Range("J134").Value = 4236.6294
$ws.Range("K134").Value = 3085668.6
$ws.Range("L134").Value = 12709.8882
$ws.Range("M134").Value = -3083133.6
$ws.Range("N134").Value = -17779.8882

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 6946199.5
$ws.Range("I31").Value = 1039.8125
$ws.Range("J31").Value = 12502327
$ws.Range("K31").Value = 1039.8125
$ws.Range("L31").Value = 12502327
$ws.Range("M31").Value = -744.8125
$ws.Range("N31").Value = -12502917
# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 6946199.5
$ws.Range("I34").Value = 1039.8125
$ws.Range("J34").Value = 12502327
$ws.Range("K34").Value = 1039.8125
$ws.Range("L34").Value = 12502327
$ws.Range("M34").Value = -837.8125
$ws.Range("N34").Value = -12502731
# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 4071.2424
$ws.Range("I58").Value = 5362.7393
$ws.Range("J58").Value = 1100.8
$ws.Range("K58").Value = 5362.7393
$ws.Range("L58").Value = 1100.8
$ws.Range("M58").Value = -5159.7393
$ws.Range("N58").Value = -1506.8
# Row 124: Earring Awakening
$ws.Range("H124").Value = 26326
$ws.Range("J124").Value = 26326
$ws.Range("L124").Value = 26326
$ws.Range("N124").Value = -31236
# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 1003062.5
$ws.Range("I132").Value = 2023.5518
$ws.Range("J132").Value = 4631828.5
$ws.Range("K132").Value = 6070.6554
$ws.Range("L132").Value = 13895485.5
$ws.Range("M132").Value = -3540.6554
$ws.Range("N132").Value = -13900545.5
# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 1855.4166
$ws.Range("I134").Value = 1882.381
$ws.Range("J134").Value = 1792.5
$ws.Range("K134").Value = 5647.143
$ws.Range("L134").Value = 5377.5
$ws.Range("M134").Value = -3112.143
$ws.Range("N134").Value = -10447.5
# Row 136: Turali Quality
$ws.Range("H136").Value = 4071.2424
$ws.Range("I136").Value = 5362.7393
$ws.Range("J136").Value = 1100.8
$ws.Range("K136").Value = 16088.2179
$ws.Range("L136").Value = 3302.4
$ws.Range("M136").Value = -13538.2179
$ws.Range("N136").Value = -8402.4

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# Row 3: Trout Fishing in Limsa
$ws.Range("H3").Value = 3766
$ws.Range("I3").Value = 1776.6666
$ws.Range("J3").Value = 6750
$ws.Range("K3").Value = 5329.9998
$ws.Range("L3").Value = 20250
$ws.Range("M3").Value = -5217.9998
$ws.Range("N3").Value = -20474
# Row 17: Chew the Fat
$ws.Range("H17").Value = 273.75
$ws.Range("I17").Value = 231.66667
$ws.Range("K17").Value = 695.00001
$ws.Range("M17").Value = -526.00001
# Row 56: Culture Club
$ws.Range("H56").Value = 5415.385
$ws.Range("I56").Value = 5415.385
$ws.Range("K56").Value = 5415.385
$ws.Range("M56").Value = -4885.385
# Row 133: Friends Are Food
$ws.Range("H133").Value = 4739.885
$ws.Range("I133").Value = 2304.111
$ws.Range("J133").Value = 6029.4116
$ws.Range("K133").Value = 6912.333
$ws.Range("L133").Value = 18088.2348
$ws.Range("M133").Value = -1852.333
$ws.Range("N133").Value = -28208.2348
# Row 137: Creative Chocolate
$ws.Range("H137").Value = 45115588
$ws.Range("I137").Value = 2732
$ws.Range("K137").Value = 8196
$ws.Range("M137").Value = -3096

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 57: Gold Is So Last Year
$ws.Range("H57").Value = 16427.092
$ws.Range("J57").Value = 17919.8
$ws.Range("L57").Value = 17919.8
$ws.Range("N57").Value = -19559.8
# Row 70: Sky Is the Limit
$ws.Range("H70").Value = 4645.8535
$ws.Range("I70").Value = 4430.5557
$ws.Range("K70").Value = 4430.5557
$ws.Range("M70").Value = -4160.5557
# Row 73: Hulls of Broken Dreams (L)
$ws.Range("H73").Value = 4645.8535
$ws.Range("I73").Value = 4430.5557
$ws.Range("K73").Value = 4430.5557
$ws.Range("M73").Value = -3494.5557
# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 130306.875
$ws.Range("I122").Value = 146700
$ws.Range("J122").Value = 15555
$ws.Range("K122").Value = 440100
$ws.Range("L122").Value = 46665
$ws.Range("M122").Value = -437650
$ws.Range("N122").Value = -51565
# Row 132: On Board for Lar
$ws.Range("H132").Value = 1726551.1
$ws.Range("I132").Value = 2010.9
$ws.Range("J132").Value = 5558862.5
$ws.Range("K132").Value = 6032.700000000001
$ws.Range("L132").Value = 16676587.5
$ws.Range("M132").Value = -3502.700000000001
$ws.Range("N132").Value = -16681647.5

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40: Best Served Toad
$ws.Range("H40").Value = 1836.8387
$ws.Range("I40").Value = 1740.826
$ws.Range("J40").Value = 2112.875
$ws.Range("K40").Value = 1740.826
$ws.Range("L40").Value = 2112.875
$ws.Range("M40").Value = -1604.826
$ws.Range("N40").Value = -2384.875
# Row 93: Hide to Go Seek
$ws.Range("H93").Value = 64864.184
$ws.Range("I93").Value = 1437.75
$ws.Range("J93").Value = 234001.33
$ws.Range("K93").Value = 1437.75
$ws.Range("L93").Value = 234001.33
$ws.Range("M93").Value = -189.75
$ws.Range("N93").Value = -236497.33

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 38: By the Seat of the Pants
$ws.Range("H38").Value = 7000
$ws.Range("J38").Value = 7000
$ws.Range("L38").Value = 7000
$ws.Range("N38").Value = -7946
# Row 132: Comfy Cabins
$ws.Range("H132").Value = 3014.6052
$ws.Range("I132").Value = 4377.1577
$ws.Range("J132").Value = 1652.0526
$ws.Range("K132").Value = 13131.4731
$ws.Range("L132").Value = 4956.1578
$ws.Range("M132").Value = -10601.4731
$ws.Range("N132").Value = -10016.1578
# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 2772.7273
$ws.Range("I136").Value = 2840.423
$ws.Range("J136").Value = 2674.9443
$ws.Range("K136").Value = 8521.269
$ws.Range("L136").Value = 8024.8329
$ws.Range("M136").Value = -5971.269
$ws.Range("N136").Value = -13124.8329

